$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.978.94'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.981.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.93%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9890'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -1.18%  '

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7685'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +63.12%  '

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '253.65'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9893'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3300'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +13.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.89'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +16.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06924'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8575'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +18.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '102.58'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07991'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.961.20'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.430'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.39'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.949.11'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.94'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007950'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +6.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.680'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.210.28'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.44%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9883'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.721'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.634'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.26'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1424'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +46.80%  '

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.70'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.193'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +15.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.567'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.361'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.533'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.333'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05101'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.225'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7412'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.693'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01985'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.923'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.605'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.17%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4740'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +11.22%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '77.95'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.078'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8505'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.03%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.09'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.75%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9906'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.931'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.584'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.10%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4278'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.60%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.20'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '943.63'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.72%  '
